# Updated cryptos list - applies Price (D) and Volume(1h) (E) column updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.543.53"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.536"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.97"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.01%  "
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.50"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.646.51"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("E16").Value = "  -2.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.558.91"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.91"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.34%  "
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0488"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.427.35"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.572"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.897"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +14.46%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.781.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0988"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.35%  "
